$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number (e.g. "1.003")
# must be forced to Text format first, otherwise Excel auto-converts the
# assigned string into a numeric value (losing e.g. trailing zeros such as
# "46.00" -> 46, or "0.07830" -> 0.0783). Prices that use "." as a thousands
# separator (e.g. "27.361.80") are never parsed as numbers, so they need no
# special handling.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '27.361.80'
$ws.Range("E2").Value = '  -1.15%  '
$ws.Range("D3").Value = '1.815.19'
$ws.Range("E3").Value = '  -3.13%  '
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -1.14%  '
$ws.Range("D5").Value = '330.80'
$ws.Range("E5").Value = '  -1.71%  '
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  -0.96%  '
$ws.Range("D7").Value = '0.4557'
$ws.Range("E7").Value = '  -2.12%  '
$ws.Range("D8").Value = '0.3802'
$ws.Range("E8").Value = '  -3.54%  '
$ws.Range("D9").Value = '46.00'
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("D10").Value = '0.07830'
$ws.Range("E10").Value = '  -2.03%  '
$ws.Range("D11").Value = '0.9575'
$ws.Range("E11").Value = '  -4.97%  '
$ws.Range("D12").Value = '20.92'
$ws.Range("E12").Value = '  -4.35%  '
$ws.Range("D13").Value = '5.836'
$ws.Range("E13").Value = '  -2.40%  '
$ws.Range("D14").Value = '1.815.48'
$ws.Range("E14").Value = '  -3.73%  '
$ws.Range("D15").Value = '7.045'
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '89.08'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '0.06581'
$ws.Range("D19").Value = '0.00001017'
$ws.Range("E19").Value = '  -2.73%  '
$ws.Range("D20").Value = '17.08'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("E21").Value = '  -0.84%  '
$ws.Range("D22").Value = '27.339.89'
$ws.Range("E22").Value = '  -1.25%  '
$ws.Range("D23").Value = '5.274'
$ws.Range("E23").Value = '  -3.59%  '
$ws.Range("D24").Value = '10.78'
$ws.Range("E24").Value = '  -1.56%  '
$ws.Range("E25").Value = '  -1.83%  '
$ws.Range("D26").Value = '2.053.67'
$ws.Range("E26").Value = '  -2.52%  '
$ws.Range("D27").Value = '155.74'
$ws.Range("E27").Value = '  -2.15%  '
$ws.Range("D28").Value = '19.26'
$ws.Range("E28").Value = '  -2.47%  '
$ws.Range("D29").Value = '2.038'
$ws.Range("E29").Value = '  -5.25%  '
$ws.Range("D30").Value = '5.228'
$ws.Range("E30").Value = '  -4.63%  '
$ws.Range("D31").Value = '117.45'
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").Value = '0.9295'
$ws.Range("E33").Value = '  -5.07%  '
$ws.Range("D34").Value = '3.565'
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("D35").Value = '5.203'
$ws.Range("E35").Value = '  -2.18%  '
$ws.Range("D36").Value = '1.308'
$ws.Range("E36").Value = '  -2.92%  '
$ws.Range("D37").Value = '0.05896'
$ws.Range("E37").Value = '  -2.74%  '
$ws.Range("D38").Value = '0.02170'
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("D39").Value = '1.003'
$ws.Range("E39").Value = '  -0.98%  '
$ws.Range("D40").Value = '8.047'
$ws.Range("E40").Value = '  -3.35%  '
$ws.Range("D41").Value = '1.135'
$ws.Range("E41").Value = '  -5.36%  '
$ws.Range("D42").Value = '0.5721'
$ws.Range("E42").Value = '  -4.03%  '
$ws.Range("E43").Value = '  -3.94%  '
$ws.Range("D44").Value = '9.888'
$ws.Range("E44").Value = '  -4.53%  '
$ws.Range("D45").Value = '1.281'
$ws.Range("E45").Value = '  +2.53%  '
$ws.Range("E46").Value = '  -3.43%  '
$ws.Range("D47").Value = '0.5367'
$ws.Range("E47").Value = '  -4.80%  '
$ws.Range("D48").Value = '1.861'
$ws.Range("E48").Value = '  -3.74%  '
$ws.Range("D49").Value = '0.06564'
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").Value = '109.46'
$ws.Range("E50").Value = '  -2.25%  '
$ws.Range("E51").Value = '  -32.39%  '
